$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 302; this pushes the existing rows
# 302-340 down to 303-341 (dimension grows from T340 to T341).
$ws.Rows(302).Insert()

# Populate the newly inserted row 302 with the new price-observation
# record (same market/region/product metadata as its neighbours).
$ws.Range("A302").Value = 10
$ws.Range("B302").Value = "Vega Modelo de Temuco"
$ws.Range("C302").Value = "La Araucanía"
$ws.Range("D302").Value = 44474
$ws.Range("E302").Value = 9
$ws.Range("F302").Value = "Fruta"
$ws.Range("G302").Value = 100108
$ws.Range("H302").Value = "Tropicales y subtropicales"
$ws.Range("I302").Value = 100108006
$ws.Range("J302").Value = "Plátano"
$ws.Range("K302").Value = "Sin especificar"
$ws.Range("L302").Value = "Pintón"
$ws.Range("M302").Value = 300
$ws.Range("N302").Value = 21000
$ws.Range("O302").Value = 21000
$ws.Range("P302").Value = 21000
$ws.Range("Q302").Value = '$/caja 20 kilos'
$ws.Range("R302").Value = "Ecuador"
$ws.Range("S302").Value = 1050
$ws.Range("T302").Value = 20
